$d = $word.ActiveDocument

# 1) "- Opcja gry dla jednego gracza / tryb test" + "owania map" (split across a
#    _GoBack bookmark) become a single merged run reading
#    "- Opcja gry dla jednego gracza / tryb testowania map" with the stray
#    bookmark removed. Searching for the full already-visible text and
#    "replacing" it with itself collapses the two runs (and the bookmark
#    sitting between them) into one clean run.
$d.Content.Find.Execute(
    "- Opcja gry dla jednego gracza / tryb testowania map", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "- Opcja gry dla jednego gracza / tryb testowania map", 2) | Out-Null

# 2) "Visual Studio 2022" -> "Visual Studio Code"
$d.Content.Find.Execute(
    "Visual Studio 2022", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Visual Studio Code", 2) | Out-Null

Write-Host "Edits applied."
